# Update cryptos list (Price / Volume(1h) columns) with refreshed quotes.
# Cells whose new text would otherwise be auto-parsed by Excel as a number
# (losing trailing zeros / decimal formatting) are written with a leading
# apostrophe so Excel stores them as plain text, matching the source data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.672.11"
$ws.Range("E2").Value = "  -0.94%  "

$ws.Range("D3").Value = "3.787.21"
$ws.Range("E3").Value = "  +1.17%  "

$ws.Range("E4").Value = "  +0.00%  "

$ws.Range("D5").Value = "'595.43"
$ws.Range("E5").Value = "  +0.33%  "

$ws.Range("D6").Value = "'167.02"
$ws.Range("E6").Value = "  +0.55%  "

$ws.Range("D7").Value = "3.773.36"
$ws.Range("E7").Value = "  +0.81%  "

$ws.Range("E8").Value = "  +0.05%  "

$ws.Range("E9").Value = "  +0.41%  "

$ws.Range("E10").Value = "  +0.00%  "

$ws.Range("D11").Value = "'6.30"
$ws.Range("E11").Value = "  -2.19%  "

$ws.Range("E12").Value = "  +0.05%  "

$ws.Range("E13").Value = "  -2.15%  "

$ws.Range("D14").Value = "'36.00"
$ws.Range("E14").Value = "  -0.52%  "

$ws.Range("D15").Value = "4.423.03"
$ws.Range("E15").Value = "  +1.22%  "

$ws.Range("D16").Value = "3.802.28"
$ws.Range("E16").Value = "  +1.58%  "

$ws.Range("D17").Value = "'18.60"
$ws.Range("E17").Value = "  +4.41%  "

$ws.Range("D18").Value = "67.654.68"
$ws.Range("E18").Value = "  -0.97%  "

$ws.Range("D19").Value = "'7.02"
$ws.Range("E19").Value = "  +0.48%  "

$ws.Range("D20").Value = "'0.112"
$ws.Range("E20").Value = "  +0.03%  "

$ws.Range("D21").Value = "'10.03"
$ws.Range("E21").Value = "  -5.96%  "

$ws.Range("D22").Value = "'459.79"
$ws.Range("E22").Value = "  -1.40%  "

$ws.Range("D23").Value = "'0.696"
$ws.Range("E23").Value = "  +0.13%  "

$ws.Range("D24").Value = "'0.0000155"
$ws.Range("E24").Value = "  +6.74%  "

$ws.Range("D25").Value = "'83.41"
$ws.Range("E25").Value = "  -0.58%  "

$ws.Range("E26").Value = "  +0.97%  "

$ws.Range("E27").Value = "  -2.68%  "

$ws.Range("E28").Value = "  +0.15%  "

$ws.Range("D29").Value = "'10.00"
$ws.Range("E29").Value = "  -0.24%  "

$ws.Range("D30").Value = "3.933.73"
$ws.Range("E30").Value = "  +1.15%  "

$ws.Range("D31").Value = "'2.77"
$ws.Range("E31").Value = "  +0.37%  "

$ws.Range("D32").Value = "'2.23"
$ws.Range("E32").Value = "  +3.32%  "

$ws.Range("D33").Value = "'7.20"
$ws.Range("E33").Value = "  -1.28%  "

$ws.Range("D34").Value = "'29.65"
$ws.Range("E34").Value = "  -0.49%  "

$ws.Range("B35").Value = "Aptos"
$ws.Range("C35").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D35").Value = "'9.09"
$ws.Range("E35").Value = "  -0.67%  "

$ws.Range("B36").Value = "Binance-PegBSC-USD"
$ws.Range("C36").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D36").Value = "'0.998"
$ws.Range("E36").Value = "  -0.21%  "

$ws.Range("E37").Value = "  -0.33%  "

$ws.Range("E38").Value = "  -0.65%  "

$ws.Range("E39").Value = "  -0.38%  "

$ws.Range("D40").Value = "'0.995"
$ws.Range("E40").Value = "  -0.36%  "

$ws.Range("E41").Value = "  -0.10%  "

$ws.Range("D42").Value = "'1.00"
$ws.Range("E42").Value = "  +0.00%  "

$ws.Range("D44").Value = "'45.21"
$ws.Range("E44").Value = "  +3.37%  "

$ws.Range("D45").Value = "'48.16"
$ws.Range("E45").Value = "  +3.25%  "

$ws.Range("E46").Value = "  -0.63%  "

$ws.Range("D47").Value = "'149.83"
$ws.Range("E47").Value = "  +3.83%  "

$ws.Range("D49").Value = "'393.89"
$ws.Range("E49").Value = "  +0.79%  "

$ws.Range("D50").Value = "'26.78"
$ws.Range("E50").Value = "  +5.81%  "

$ws.Range("E51").Value = "  -4.90%  "
